# Generate Report for Handoff
#
# A new handoff xliff was generated for 498504b3-d390-4808-a139-79672ffeb6c6:
#   - zh-cn "Latest Handoff Datetime" moves from 2016-08-31 08:51:13 to 2016-08-31 08:51:34
#   - Overview "Latest HO Xliff Generate Date" (row for the same file) moves from
#     2016-08-31 08:51:19 to 2016-08-31 08:51:39
#
# All other apparent differences in the underlying shared-string table are just
# bookkeeping side effects of inserting the new strings, not content changes.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Overview sheet: row 6 is 498504b3-d390-4808-a139-79672ffeb6c6.md, column G = "Latest HO Xliff Generate Date"
$wsOverview.Cells.Item(6, 7).Value = "2016-08-31 08:51:39"

# zh-cn sheet: row 6 is 498504b3-d390-4808-a139-79672ffeb6c6.md, column H = "Latest Handoff Datetime"
$wsZhCn.Cells.Item(6, 8).Value = "2016-08-31 08:51:34"
